$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the "Förändrad" (changed) date column C for all data rows (2..302)
#    from 45181 (2023-09-29) to 45182 (2023-09-30).
$ws.Range("C2:C302").Value = 45182

# 2. Row 302 gains an explicit row height (ht="15" customHeight="1").
$ws.Rows.Item(302).RowHeight = 15

# 3. Append three new data rows (303, 304, 305) at the bottom of the table.
#    Rows 303 and 304 get an explicit custom row height like row 302;
#    row 305 (the new last row) keeps the default (no explicit height).
$newRows = @(
    @{ Row = 303; Beteckning = "A 42733-2023"; Datum = 45181; Forandrad = 45182; Area = 2.7 },
    @{ Row = 304; Beteckning = "A 42673-2023"; Datum = 45181; Forandrad = 45182; Area = 0.3 },
    @{ Row = 305; Beteckning = "A 42753-2023"; Datum = 45181; Forandrad = 45182; Area = 0.5 }
)

foreach ($r in $newRows) {
    $row = $r.Row

    $ws.Cells.Item($row, 1).Value = $r.Beteckning

    $ws.Cells.Item($row, 2).Value = $r.Datum
    $ws.Cells.Item($row, 2).NumberFormat = "YYYY-MM-DD"

    $ws.Cells.Item($row, 3).Value = $r.Forandrad
    $ws.Cells.Item($row, 3).NumberFormat = "YYYY-MM-DD"

    $ws.Cells.Item($row, 4).Value = "ÖSTERGÖTLANDS LÄN"
    $ws.Cells.Item($row, 5).Value = "SÖDERKÖPING"

    $ws.Cells.Item($row, 7).Value = $r.Area
    $ws.Cells.Item($row, 8).Value = 0
    $ws.Cells.Item($row, 9).Value = 0
    $ws.Cells.Item($row, 10).Value = 0
    $ws.Cells.Item($row, 11).Value = 0
    $ws.Cells.Item($row, 12).Value = 0
    $ws.Cells.Item($row, 13).Value = 0
    $ws.Cells.Item($row, 14).Value = 0
    $ws.Cells.Item($row, 15).Value = 0
    $ws.Cells.Item($row, 16).Value = 0
    $ws.Cells.Item($row, 17).Value = 0

    $ws.Cells.Item($row, 18).Value = ""
    $ws.Cells.Item($row, 18).WrapText = $true
}

$ws.Rows.Item(303).RowHeight = 15
$ws.Rows.Item(304).RowHeight = 15
